$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("O2").Value = 1.73
$ws.Range("P2").Value = 2
$ws.Range("V2").Value = 1.47

# Row 3 updates
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 19
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 5.5
$ws.Range("Q3").Value = 1.5
$ws.Range("R3").Value = 2.5
